$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Wine row (currently row 42) gets its Description filled in first.
#    This is the first new piece of text typed, so it becomes shared string 116.
$ws.Range("B42").Value = "An alcoholic drink derived from grapes. Quite common at the dinner table and an important addition at religious ceremonies"

# 2) Insert a new row for GRAPES right before row 20 (HOPS), pushing
#    everything from row 20 down by one.
$ws.Rows.Item(20).Insert()

# Copy the format of the row above (row 19) onto the new row 20 so the
# cell styles (borders/fonts) match the rest of the table, then fill in
# the new resource's data.
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = "GRAPES"
$ws.Range("B20").Value = "A round fruit grown on vines known for its sweetness. Mostly used to make wines. Could be used to make jams and candied grapes."
$ws.Range("C20").ClearContents()
$ws.Rows.Item(20).RowHeight = 30

# 3) Wine row is now row 43 (shifted by the insert above). Fill in its
#    Price column and bump the row height to match the two other updated
#    cells.
$ws.Range("C43").Value = "s"
$ws.Rows.Item(43).RowHeight = 30

# Leave the selection where the author ended up (Price column of the
# Wine row).
[void]$ws.Range("C43").Select()
